$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.0003714022599530242
$ws.Range("C2").Value = 0.0001537489499301437
$ws.Range("D2").Value = 157.8057217802531
$ws.Range("E2").Value = 246.9852506941017
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 404.7914976255647
